$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text to uppercase labels
$ws.Range("A1").Value = "AKAN"
$ws.Range("B1").Value = "ENGLISH"

# Remove the stray formatted-but-empty C1 cell
$ws.Range("C1").Clear()

# Header cells should not wrap (matches the remaining style definition)
$ws.Range("A1:B1").WrapText = $false

# Update the saved selection to the whole header row, matching the author's last selection
[void]$ws.Range("A1:XFD1").Select()
